# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# OFF sheet - Home row (row 2) target depth totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 202
$wsOff.Range("C2").Value = 150
$wsOff.Range("D2").Value = 59
$wsOff.Range("E2").Value = 29
$wsOff.Range("G2").Value = 3

# DEF sheet - Home row (row 2) target depth totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 196
$wsDef.Range("C2").Value = 130
$wsDef.Range("D2").Value = 45
$wsDef.Range("E2").Value = 16
$wsDef.Range("G2").Value = 7
